# Update crypto price/volume data per GitHub Actions scrape run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.043.78"
$ws.Range("E2").Value = "  -0.58%  "

# Row 3
$ws.Range("D3").Value = "1.827.08"
$ws.Range("E3").Value = "  -0.50%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9991"
$ws.Range("D4").Style = "Normal"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "241.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.09%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6371"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.46%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("E8").Value = "  +6.51%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2937"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.38%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07342"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.48%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "22.81"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.55%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07667"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.68%  "

# Row 13
$ws.Range("D13").Value = "1.824.51"
$ws.Range("E13").Value = "  +0.43%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.987"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.03%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6635"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.83%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.98"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.59%  "

# Row 17
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008685"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.05%  "

# Row 18
$ws.Range("B18").Value = "Uniswap"
$ws.Range("C18").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.052"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.98%  "

# Row 19
$ws.Range("D19").Value = "28.876.61"
$ws.Range("E19").Value = "  -1.01%  "

# Row 20
$ws.Range("D20").Value = "2.074.90"
$ws.Range("E20").Value = "  +0.42%  "

# Row 21
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "224.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.49%  "

# Row 22
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.61%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.00%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.120"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.21%  "

# Row 25
$ws.Range("E25").Value = "  +0.05%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.63%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.473"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.96%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1372"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.67%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "17.86"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.52%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.505"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.26%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.092"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.59%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.026"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.15%  "

# Row 33
$ws.Range("E33").Value = "  +1.70%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05290"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.40%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.834"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.89%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7368"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.20%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.154"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.05%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.652"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.96%  "

# Row 39
$ws.Range("D39").Value = "1.295.14"
$ws.Range("E39").Value = "  -0.10%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.749"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.06%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01782"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.83%  "

# Row 42
$ws.Range("E42").Value = "  +5.50%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8988"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.18%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9997"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.53%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "102.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.11%  "

# Row 46
$ws.Range("D46").Value = "1.974.89"
$ws.Range("E46").Value = "  +0.20%  "

# Row 47
$ws.Range("E47").Value = "  -0.46%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "64.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.38%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00000000119"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.51%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.726"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.46%  "

# Row 51
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05800"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.23%  "
